$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 368.5
$ws.Range("I6").Value = 368.5
$ws.Range("K6").Value = 1105.5
$ws.Range("M6").Value = -993.5

$ws.Range("H53").Value = 260
$ws.Range("I53").Value = 225
$ws.Range("J53").Value = 321.25
$ws.Range("K53").Value = 225
$ws.Range("L53").Value = 321.25
$ws.Range("M53").Value = 412
$ws.Range("N53").Value = -1595.25

$ws.Range("H62").Value = 5573.391
$ws.Range("I62").Value = 4604.8887
$ws.Range("K62").Value = 4604.8887
$ws.Range("M62").Value = -3980.8887

$ws.Range("H65").Value = 5573.391
$ws.Range("I65").Value = 4604.8887
$ws.Range("K65").Value = 23024.4435
$ws.Range("M65").Value = -19904.4435

$ws.Range("H70").Value = 3289.739
$ws.Range("I70").Value = 1826.2858
$ws.Range("K70").Value = 5478.857400000001
$ws.Range("M70").Value = -5208.857400000001

$ws.Range("H73").Value = 3289.739
$ws.Range("I73").Value = 1826.2858
$ws.Range("K73").Value = 5478.857400000001
$ws.Range("M73").Value = -4542.857400000001

$ws.Range("H76").Value = 6432.222
$ws.Range("I76").Value = 5912.857
$ws.Range("K76").Value = 5912.857
$ws.Range("M76").Value = -5597.857

$ws.Range("H79").Value = 6432.222
$ws.Range("I79").Value = 5912.857
$ws.Range("K79").Value = 5912.857
$ws.Range("M79").Value = -4820.857

$ws.Range("H98").Value = 807.7778
$ws.Range("I98").Value = 767.0833
$ws.Range("K98").Value = 767.0833
$ws.Range("M98").Value = 730.9167

$ws.Range("H122").Value = 807.7778
$ws.Range("I122").Value = 767.0833
$ws.Range("K122").Value = 2301.2499
$ws.Range("M122").Value = 148.7501000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 16749.5
$ws.Range("I8").Value = 12333
$ws.Range("J8").Value = 29999
$ws.Range("K8").Value = 12333
$ws.Range("L8").Value = 29999
$ws.Range("M8").Value = -12189
$ws.Range("N8").Value = -30287

$ws.Range("H26").Value = 2962.5
$ws.Range("I26").Value = 2962.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 2962.5
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -2632.5
$ws.Range("N26").ClearContents()

$ws.Range("H63").Value = 1869.3
$ws.Range("I63").Value = 1682
$ws.Range("K63").Value = 1682
$ws.Range("M63").Value = -996

$ws.Range("H66").Value = 1869.3
$ws.Range("I66").Value = 1682
$ws.Range("K66").Value = 8410
$ws.Range("M66").Value = -4978

$ws.Range("H122").Value = 2735.8076
$ws.Range("I122").Value = 2363.5715
$ws.Range("J122").Value = 4299.2
$ws.Range("K122").Value = 7090.7145
$ws.Range("L122").Value = 12897.6
$ws.Range("M122").Value = -4640.7145
$ws.Range("N122").Value = -17797.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 77499
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H94").Value = 9989.5
$ws.Range("I94").Value = 9989.5
$ws.Range("K94").Value = 9989.5
$ws.Range("M94").Value = -9538.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 16842.857
$ws.Range("I42").Value = 9300
$ws.Range("J42").Value = 22500
$ws.Range("K42").Value = 9300
$ws.Range("L42").Value = 22500
$ws.Range("M42").Value = -8707
$ws.Range("N42").Value = -23686

$ws.Range("H134").Value = 4023.6072
$ws.Range("I134").Value = 3598.48
$ws.Range("J134").Value = 7566.3335
$ws.Range("K134").Value = 10795.44
$ws.Range("L134").Value = 22699.0005
$ws.Range("M134").Value = -8260.440000000001
$ws.Range("N134").Value = -27769.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1078.4
$ws.Range("I14").Value = 1078.4
$ws.Range("K14").Value = 3235.2
$ws.Range("M14").Value = -3062.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 996.9091
$ws.Range("I2").Value = 1144.5625
$ws.Range("K2").Value = 1144.5625
$ws.Range("M2").Value = -1031.5625

$ws.Range("H19").Value = 9995
$ws.Range("I19").Value = 9995
$ws.Range("K19").Value = 9995
$ws.Range("M19").Value = -9707

$ws.Range("H35").Value = 111555
$ws.Range("J35").Value = 111555
$ws.Range("L35").Value = 111555
$ws.Range("N35").Value = -112151

$ws.Range("H62").Value = 80000
$ws.Range("J62").Value = 80000
$ws.Range("L62").Value = 80000
$ws.Range("N62").Value = -81372

$ws.Range("H65").Value = 80000
$ws.Range("J65").Value = 80000
$ws.Range("L65").Value = 240000
$ws.Range("N65").Value = -246864

$ws.Range("H132").Value = 5604.8066
$ws.Range("I132").Value = 4236.278
$ws.Range("J132").Value = 7499.6924
$ws.Range("K132").Value = 12708.834
$ws.Range("L132").Value = 22499.0772
$ws.Range("M132").Value = -10178.834
$ws.Range("N132").Value = -27559.0772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1367
$ws.Range("I4").Value = 1367
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1367
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -1254
$ws.Range("N4").ClearContents()

$ws.Range("H28").Value = 1367
$ws.Range("I28").Value = 1367
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1367
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -1135
$ws.Range("N28").ClearContents()

$ws.Range("H37").Value = 1367
$ws.Range("I37").Value = 1367
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1367
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -1260
$ws.Range("N37").ClearContents()

$ws.Range("H46").Value = 3397.0625
$ws.Range("I46").Value = 2023
$ws.Range("J46").Value = 3855.0833
$ws.Range("K46").Value = 2023
$ws.Range("L46").Value = 3855.0833
$ws.Range("M46").Value = -1835
$ws.Range("N46").Value = -4231.0833

$ws.Range("H122").Value = 9726.799999999999
$ws.Range("I122").Value = 8939.700000000001
$ws.Range("K122").Value = 26819.1
$ws.Range("M122").Value = -24369.1

$ws.Range("H132").Value = 3559.4783
$ws.Range("I132").Value = 2263.9
$ws.Range("K132").Value = 6791.700000000001
$ws.Range("M132").Value = -4261.700000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
